$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2427.5334
$ws.Cells.Item(32, 9).Value = 2107
$ws.Cells.Item(32, 10).Value = 2641.2222
$ws.Cells.Item(32, 11).Value = 2107
$ws.Cells.Item(32, 12).Value = 2641.2222
$ws.Cells.Item(32, 13).Value = -1781
$ws.Cells.Item(32, 14).Value = -3293.2222
$ws.Cells.Item(80, 8).Value = 436799.88
$ws.Cells.Item(80, 10).Value = 836491.7
$ws.Cells.Item(80, 12).Value = 2509475.1
$ws.Cells.Item(80, 14).Value = -2511471.1
$ws.Cells.Item(83, 8).Value = 436799.88
$ws.Cells.Item(83, 10).Value = 836491.7
$ws.Cells.Item(83, 12).Value = 7528425.3
$ws.Cells.Item(83, 14).Value = -7538409.3
$ws.Cells.Item(106, 8).Value = 30319.842
$ws.Cells.Item(106, 9).Value = 7278.8
$ws.Cells.Item(106, 11).Value = 7278.8
$ws.Cells.Item(106, 13).Value = -6647.8
$ws.Cells.Item(113, 8).Value = 4441.2334
$ws.Cells.Item(113, 9).Value = 3901.625
$ws.Cells.Item(113, 10).Value = 6599.6665
$ws.Cells.Item(113, 11).Value = 3901.625
$ws.Cells.Item(113, 12).Value = 6599.6665
$ws.Cells.Item(113, 13).Value = -647.625
$ws.Cells.Item(113, 14).Value = -13107.6665

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1049.6666
$ws.Cells.Item(5, 9).Value = 1066
$ws.Cells.Item(5, 10).Value = 1038.7778
$ws.Cells.Item(5, 11).Value = 1066
$ws.Cells.Item(5, 12).Value = 1038.7778
$ws.Cells.Item(5, 13).Value = -954
$ws.Cells.Item(5, 14).Value = -1262.7778
$ws.Cells.Item(32, 8).Value = 17206.656
$ws.Cells.Item(32, 9).Value = 15883.207
$ws.Cells.Item(32, 10).Value = 30000
$ws.Cells.Item(32, 11).Value = 15883.207
$ws.Cells.Item(32, 12).Value = 30000
$ws.Cells.Item(32, 13).Value = -15596.207
$ws.Cells.Item(32, 14).Value = -30574
$ws.Cells.Item(61, 8).Value = 11219.4375
$ws.Cells.Item(61, 9).Value = 3601.3635
$ws.Cells.Item(61, 11).Value = 3601.3635
$ws.Cells.Item(61, 13).Value = -3389.3635
$ws.Cells.Item(102, 8).Value = 3353.9412
$ws.Cells.Item(102, 9).Value = 1063.625
$ws.Cells.Item(102, 11).Value = 1063.625
$ws.Cells.Item(102, 13).Value = 558.375
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 3513.6758
$ws.Cells.Item(132, 9).Value = 3513.6758
$ws.Cells.Item(132, 11).Value = 10541.0274
$ws.Cells.Item(132, 13).Value = -8011.027399999999
$ws.Cells.Item(136, 8).Value = 11219.4375
$ws.Cells.Item(136, 9).Value = 3601.3635
$ws.Cells.Item(136, 11).Value = 10804.0905
$ws.Cells.Item(136, 13).Value = -8254.0905

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1049.6666
$ws.Cells.Item(4, 9).Value = 1066
$ws.Cells.Item(4, 10).Value = 1038.7778
$ws.Cells.Item(4, 11).Value = 1066
$ws.Cells.Item(4, 12).Value = 1038.7778
$ws.Cells.Item(4, 13).Value = -951
$ws.Cells.Item(4, 14).Value = -1268.7778

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 39178
$ws.Cells.Item(16, 9).Value = 2479.5
$ws.Cells.Item(16, 11).Value = 2479.5
$ws.Cells.Item(16, 13).Value = -2192.5
$ws.Cells.Item(41, 8).Value = 34994
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 14).ClearContents()
$ws.Cells.Item(59, 8).Value = 74996
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).ClearContents()
$ws.Cells.Item(60, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 40000
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 40000
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 3248.6924
$ws.Cells.Item(86, 10).Value = 3900.7222
$ws.Cells.Item(86, 12).Value = 3900.7222
$ws.Cells.Item(86, 14).Value = -6146.7222
$ws.Cells.Item(89, 8).Value = 3248.6924
$ws.Cells.Item(89, 10).Value = 3900.7222
$ws.Cells.Item(89, 12).Value = 19503.611
$ws.Cells.Item(89, 14).Value = -30735.611
$ws.Cells.Item(99, 8).Value = 10899.325
$ws.Cells.Item(99, 9).Value = 7598.1
$ws.Cells.Item(99, 10).Value = 11899.697
$ws.Cells.Item(99, 11).Value = 7598.1
$ws.Cells.Item(99, 12).Value = 11899.697
$ws.Cells.Item(99, 13).Value = -6100.1
$ws.Cells.Item(99, 14).Value = -14895.697
$ws.Cells.Item(113, 8).Value = 39178
$ws.Cells.Item(113, 9).Value = 2479.5
$ws.Cells.Item(113, 11).Value = 2479.5
$ws.Cells.Item(113, 13).Value = -309.5
$ws.Cells.Item(126, 8).Value = 10899.325
$ws.Cells.Item(126, 9).Value = 7598.1
$ws.Cells.Item(126, 10).Value = 11899.697
$ws.Cells.Item(126, 11).Value = 22794.3
$ws.Cells.Item(126, 12).Value = 35699.091
$ws.Cells.Item(126, 13).Value = -20324.3
$ws.Cells.Item(126, 14).Value = -40639.091

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 1801
$ws.Cells.Item(12, 9).Value = 27.25
$ws.Cells.Item(12, 11).Value = 81.75
$ws.Cells.Item(12, 13).Value = 91.25
$ws.Cells.Item(26, 8).Value = 81.75
$ws.Cells.Item(26, 9).Value = 91.22221999999999
$ws.Cells.Item(26, 10).Value = 69.57143000000001
$ws.Cells.Item(26, 11).Value = 273.66666
$ws.Cells.Item(26, 12).Value = 208.71429
$ws.Cells.Item(26, 13).Value = 14.33334000000002
$ws.Cells.Item(26, 14).Value = -784.71429
$ws.Cells.Item(129, 8).Value = 2986.0667
$ws.Cells.Item(129, 10).Value = 3826.5
$ws.Cells.Item(129, 12).Value = 11479.5
$ws.Cells.Item(129, 14).Value = -21479.5
$ws.Cells.Item(131, 8).Value = 5039.109
$ws.Cells.Item(131, 9).Value = 1576.1666
$ws.Cells.Item(131, 10).Value = 7265.2856
$ws.Cells.Item(131, 11).Value = 4728.4998
$ws.Cells.Item(131, 12).Value = 21795.8568
$ws.Cells.Item(131, 13).Value = 311.5002000000004
$ws.Cells.Item(131, 14).Value = -31875.8568

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 663.5
$ws.Cells.Item(107, 9).Value = 767.8
$ws.Cells.Item(107, 10).Value = 440
$ws.Cells.Item(107, 11).Value = 767.8
$ws.Cells.Item(107, 12).Value = 440
$ws.Cells.Item(107, 13).Value = 1152.2
$ws.Cells.Item(107, 14).Value = -4280
$ws.Cells.Item(122, 8).Value = 2359.625
$ws.Cells.Item(122, 9).Value = 2411.2856
$ws.Cells.Item(122, 11).Value = 7233.8568
$ws.Cells.Item(122, 13).Value = -4783.8568
$ws.Cells.Item(126, 8).Value = 4723.7334
$ws.Cells.Item(126, 10).Value = 6030.6
$ws.Cells.Item(126, 12).Value = 18091.8
$ws.Cells.Item(126, 14).Value = -23031.8
$ws.Cells.Item(132, 8).Value = 6626.7827
$ws.Cells.Item(132, 9).Value = 6746.1816
$ws.Cells.Item(132, 11).Value = 20238.5448
$ws.Cells.Item(132, 13).Value = -17708.5448

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5173.5713
$ws.Cells.Item(7, 9).Value = 3517.5833
$ws.Cells.Item(7, 10).Value = 15109.5
$ws.Cells.Item(7, 11).Value = 3517.5833
$ws.Cells.Item(7, 12).Value = 15109.5
$ws.Cells.Item(7, 13).Value = -3405.5833
$ws.Cells.Item(7, 14).Value = -15333.5
$ws.Cells.Item(18, 8).Value = 3872.6667
$ws.Cells.Item(18, 10).Value = 3872.6667
$ws.Cells.Item(18, 12).Value = 3872.6667
$ws.Cells.Item(18, 14).Value = -4216.6667
$ws.Cells.Item(61, 8).Value = 3110.516
$ws.Cells.Item(61, 9).Value = 2770.2693
$ws.Cells.Item(61, 11).Value = 2770.2693
$ws.Cells.Item(61, 13).Value = -2568.2693
$ws.Cells.Item(113, 8).Value = 3110.516
$ws.Cells.Item(113, 9).Value = 2770.2693
$ws.Cells.Item(113, 11).Value = 2770.2693
$ws.Cells.Item(113, 13).Value = -600.2692999999999
$ws.Cells.Item(122, 8).Value = 5417.1904
$ws.Cells.Item(122, 10).Value = 6070.8184
$ws.Cells.Item(122, 12).Value = 18212.4552
$ws.Cells.Item(122, 14).Value = -23112.4552
$ws.Cells.Item(126, 8).Value = 5173.5713
$ws.Cells.Item(126, 9).Value = 3517.5833
$ws.Cells.Item(126, 10).Value = 15109.5
$ws.Cells.Item(126, 11).Value = 10552.7499
$ws.Cells.Item(126, 12).Value = 45328.5
$ws.Cells.Item(126, 13).Value = -8082.749899999999
$ws.Cells.Item(126, 14).Value = -50268.5
$ws.Cells.Item(132, 8).Value = 7624.154
$ws.Cells.Item(132, 9).Value = 7189.161
$ws.Cells.Item(132, 10).Value = 9309.75
$ws.Cells.Item(132, 11).Value = 21567.483
$ws.Cells.Item(132, 12).Value = 27929.25
$ws.Cells.Item(132, 13).Value = -19037.483
$ws.Cells.Item(132, 14).Value = -32989.25
$ws.Cells.Item(133, 8).Value = 87608
$ws.Cells.Item(133, 10).Value = 87608
$ws.Cells.Item(133, 12).Value = 87608
$ws.Cells.Item(133, 14).Value = -92668
$ws.Cells.Item(136, 8).Value = 8787.654
$ws.Cells.Item(136, 9).Value = 8675.022999999999
$ws.Cells.Item(136, 11).Value = 26025.069
$ws.Cells.Item(136, 13).Value = -23475.069

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 84714
$ws.Cells.Item(46, 10).Value = 84714
$ws.Cells.Item(46, 12).Value = 84714
$ws.Cells.Item(46, 14).Value = -85176
$ws.Cells.Item(113, 8).Value = 915.1539
$ws.Cells.Item(113, 10).Value = 879.6
$ws.Cells.Item(113, 12).Value = 2638.8
$ws.Cells.Item(113, 14).Value = -6978.8
$ws.Cells.Item(122, 8).Value = 125233.65
$ws.Cells.Item(122, 10).Value = 100498.45
$ws.Cells.Item(122, 12).Value = 301495.35
$ws.Cells.Item(122, 14).Value = -306395.35
$ws.Cells.Item(132, 8).Value = 198019.8
$ws.Cells.Item(132, 9).Value = 285921.6
$ws.Cells.Item(132, 10).Value = 35401.45
$ws.Cells.Item(132, 11).Value = 857764.7999999999
$ws.Cells.Item(132, 12).Value = 106204.35
$ws.Cells.Item(132, 13).Value = -855234.7999999999
$ws.Cells.Item(132, 14).Value = -111264.35
$ws.Cells.Item(134, 8).Value = 84714
$ws.Cells.Item(134, 10).Value = 84714
$ws.Cells.Item(134, 12).Value = 254142
$ws.Cells.Item(134, 14).Value = -259212
$ws.Cells.Item(136, 8).Value = 7694785
$ws.Cells.Item(136, 9).Value = 11539948
$ws.Cells.Item(136, 11).Value = 34619844
$ws.Cells.Item(136, 13).Value = -34617294

Write-Host "All updates applied."